# edit.ps1 - applies "edited repo location and note" commit
#
# Summary of changes (see xml_diff):
#  Slide 3 (Content Placeholder 2):
#    - "A little overview of some terminology behind web development"
#        -> word "little" replaced with "short" (splits the run in three:
#           "A " / "short " / "overview of some terminology behind web development")
#    - "A note about JSON and maximizing data portability"
#        -> replaced entirely with "Setting up a local webserver"
#  Slide 4 (Content Placeholder 2):
#    - "http://" + "bit.ly/2mRzSc3" runs merged & URL updated to
#        "http://bit.ly/2F2RrgF"
#    - "MDN " / "(Mozilla Developer Network) ... out " / "there" runs merged
#        into a single run with the full sentence
#    - "Guide to Dynamic " run split into "Guide " / "to Dynamic "
#  (endParaRPr additions/removals scattered through the diff are incidental
#   PowerPoint-editor bookkeeping around the runs above - not independently
#   reproducible through the exposed COM surface - so they are left as-is.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: "What this workshop focuses on:"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 4: "A little overview of some terminology behind web development"
# Replace the word "little " with "short " - this splits the single run
# into three runs: "A " / "short " / "overview ... development"
$para = $body3.Paragraphs(4, 1)
$word = $body3.Characters($para.Start + 2, 7)   # "little "
$word.Text = "short "

# Paragraph 5: "A note about JSON and maximizing data portability"
# Replace with the new bullet text entirely.
$para = $body3.Paragraphs(5, 1)
$para.Text = "Setting up a local webserver"

# ---------------------------------------------------------------------
# Slide 4: "Resources for today:"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "The repository is here: http://bit.ly/2mRzSc3"
# Replace from "http://" to the end of the paragraph with the new link -
# this merges the two old runs ("http://" / "bit.ly/2mRzSc3") into one.
$para = $body4.Paragraphs(1, 1)
$tail = $body4.Characters($para.Start + 24, $para.Length - 24 - 1)
$tail.Text = "http://bit.ly/2F2RrgF"

# Paragraph 2: "MDN (Mozilla Developer Network) ... out there"
# Re-assign the identical text over the whole paragraph span so the three
# existing runs collapse into a single run.
$para = $body4.Paragraphs(2, 1)
$whole = $body4.Characters($para.Start, $para.Length - 1)
$whole.Text = $whole.Text

# Paragraph 3 (level 1): "Guide to Dynamic Javascript (Client-Side)"
# Split "Guide to Dynamic " into "Guide " / "to Dynamic ".
$para = $body4.Paragraphs(3, 1)
$word = $body4.Characters($para.Start + 6, 11)  # "to Dynamic "
$word.Text = "to Dynamic "
